# Updates for "GUERRERO FAREZ FABIAN MAURICIO" workbook
# A new client, "CONSORCIO MUELLE CENTRO DE ARTE", is inserted alphabetically
# (between "CERAMIKASA S.A.S." and "DECOGARCIA S.A.S.") on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, with all-zero figures.
# Inserting the row shifts every following row down by one, and the final
# totals/summary rows move down accordingly (their "X de N" counters on the
# VENTAS POR GRUPO sheet bump their denominator from 52 to 53).

$wb = $excel.ActiveWorkbook

$advisor = "GUERRERO FAREZ FABIAN MAURICIO"
$newClient = "CONSORCIO MUELLE CENTRO DE ARTE"
$insertRow = 15

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R, data rows 2-53, totals row 54)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push every row from 15 down onward one row lower, creating a blank row 15
$ws1.Rows.Item($insertRow).Insert()

# Populate the newly inserted row for the new client (all zero amounts)
$ws1.Cells.Item($insertRow, 1).Value = $advisor
$ws1.Cells.Item($insertRow, 2).Value = $newClient
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item($insertRow, $col).Value = 0
}

# The totals row (originally row 54, now row 55) shows "<n> de 52" counts;
# bump the denominator to 53 to reflect the extra client row.
$totalsRow1 = 55
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item($totalsRow1, $col)
    $cell.Value = $cell.Value() -replace " de 52$", " de 53"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G, data rows 2-57, totals row 58)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item($insertRow).Insert()

$ws2.Cells.Item($insertRow, 1).Value = $advisor
$ws2.Cells.Item($insertRow, 2).Value = $newClient
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item($insertRow, $col).Value = 0
}
